# Updates the cryptos list (price + volume columns) as published by the
# "Updated cryptos list" GitHub Actions workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($row, $value)
    $cell = $ws.Range("D$row")
    # Force text storage so numeric-looking strings (e.g. "5.40", "90.10")
    # are not silently coerced into floating point numbers by Excel, which
    # would drop significant trailing zeros / merge the thousands dots.
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

function Set-VolumeText {
    param($row, $value)
    $ws.Range("E$row").Value = $value
}

# Row 2 - Bitcoin
Set-PriceText 2 "26.682.26"
Set-VolumeText 2 "  -0.17%  "

# Row 3 - Ethereum
Set-PriceText 3 "1.600.73"
Set-VolumeText 3 "  +0.40%  "

# Row 4 - TetherUSD
Set-VolumeText 4 "  -0.05%  "

# Row 5 - BNB
Set-PriceText 5 "211.48"
Set-VolumeText 5 "  -0.08%  "

# Row 6 - XRP
Set-VolumeText 6 "  +1.10%  "

# Row 7 - USDC
Set-VolumeText 7 "  -0.06%  "

# Row 8 - Dogecoin
Set-VolumeText 8 "  -0.02%  "

# Row 9 - Cardano
Set-PriceText 9 "0.245"
Set-VolumeText 9 "  -1.27%  "

# Row 10 - Solana
Set-VolumeText 10 "  +0.52%  "

# Row 11 - TRON
Set-VolumeText 11 "  +0.11%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-PriceText 12 "1.824.12"
Set-VolumeText 12 "  +0.26%  "

# Row 13 - WrappedEther
Set-PriceText 13 "1.600.44"
Set-VolumeText 13 "  +0.00%  "

# Row 14 - Polkadot
Set-PriceText 14 "4.04"
Set-VolumeText 14 "  -0.31%  "

# Row 15 - Polygon
Set-VolumeText 15 "  -1.38%  "

# Row 16 - Litecoin
Set-PriceText 16 "64.88"
Set-VolumeText 16 "  +1.98%  "

# Row 17 - WrappedBTC
Set-PriceText 17 "26.664.64"
Set-VolumeText 17 "  -0.39%  "

# Row 18 - ShibaInu
Set-PriceText 18 "0.0₃0729"
Set-VolumeText 18 "  -0.10%  "

# Row 19 - BitcoinCash
Set-PriceText 19 "209.79"
Set-VolumeText 19 "  +0.45%  "

# Row 20 - Dai
Set-VolumeText 20 "  +0.06%  "

# Row 21 - Chainlink
Set-VolumeText 21 "  +1.00%  "

# Row 22 - Uniswap
Set-VolumeText 22 "  +0.35%  "

# Row 23 - Toncoin
Set-PriceText 23 "2.29"
Set-VolumeText 23 "  -2.68%  "

# Row 24 - Avalanche
Set-PriceText 24 "8.92"
Set-VolumeText 24 "  +0.64%  "

# Row 25 - Monero
Set-PriceText 25 "146.42"
Set-VolumeText 25 "  -0.08%  "

# Row 26 - BinanceUSD
Set-VolumeText 26 "  -0.12%  "

# Row 27 - Cosmos
Set-PriceText 27 "7.21"
Set-VolumeText 27 "  -3.43%  "

# Row 28 - Stellar
Set-VolumeText 28 "  +2.53%  "

# Row 29 - EthereumClassic
Set-VolumeText 29 "  -0.11%  "

# Row 30 - Hedera
Set-PriceText 30 "0.0505"
Set-VolumeText 30 "  +1.06%  "

# Row 31 - PancakeSwap
Set-VolumeText 31 "  -0.31%  "

# Row 32 - Filecoin
Set-VolumeText 32 "  -0.85%  "

# Row 33 - ImmutableX
Set-PriceText 33 "0.668"
Set-VolumeText 33 "  +0.04%  "

# Row 34 - InternetComputer(DFINITY)
Set-PriceText 34 "2.93"
Set-VolumeText 34 "  -0.41%  "

# Row 35 - Maker
Set-PriceText 35 "1.297.07"
Set-VolumeText 35 "  -1.15%  "

# Row 37 - LidoDAOToken
Set-PriceText 37 "1.48"
Set-VolumeText 37 "  -1.49%  "

# Row 38 - VeChain
Set-VolumeText 38 "  -0.87%  "

# Row 39 - ARBITRUM
Set-PriceText 39 "0.843"
Set-VolumeText 39 "  +2.91%  "

# Row 40 - PaxDollar
Set-VolumeText 40 "  -0.02%  "

# Row 41 - FraxShare
Set-PriceText 41 "5.40"
Set-VolumeText 41 "  +1.37%  "

# Rows 42/43 - MXToken and TrustWalletToken swap ranking order.
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-PriceText 42 "0.790"
Set-VolumeText 42 "  +0.19%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-PriceText 43 "2.20"
Set-VolumeText 43 "  +1.11%  "

# Row 44 - Aave
Set-PriceText 44 "63.94"
Set-VolumeText 44 "  +1.35%  "

# Row 45 - RocketPoolETH
Set-PriceText 45 "1.736.48"
Set-VolumeText 45 "  +0.25%  "

# Row 46 - WEMIXToken
Set-PriceText 46 "0.896"
Set-VolumeText 46 "  +9.09%  "

# Row 47 - Quant
Set-PriceText 47 "90.10"
Set-VolumeText 47 "  +1.19%  "

# Row 48 - RenderToken
Set-VolumeText 48 "  +0.19%  "

# Row 49 - Algorand
Set-VolumeText 49 "  +2.36%  "

# Row 50 - Cronos
Set-PriceText 50 "0.0504"
Set-VolumeText 50 "  -0.93%  "

# Row 51 - EnergySwap
Set-PriceText 51 "7.48"
Set-VolumeText 51 "  +0.47%  "
